$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) cells to be treated as text so that values like
# "0.999" or "69.057.40" are not auto-converted to numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '69.057.40'
$ws.Range('E2').Value = '  -0.98%  '
$ws.Range('D3').Value = '3.509.46'
$ws.Range('E3').Value = '  -1.35%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '569.16'
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('D6').Value = '181.97'
$ws.Range('E6').Value = '  -3.15%  '
$ws.Range('D7').Value = '3.503.57'
$ws.Range('E7').Value = '  -1.34%  '
$ws.Range('E8').Value = '  -2.65%  '
$ws.Range('E9').Value = '  -0.02%  '
$ws.Range('D10').Value = '0.188'
$ws.Range('E10').Value = '  +6.87%  '
$ws.Range('E11').Value = '  -2.81%  '
$ws.Range('D12').Value = '53.74'
$ws.Range('E12').Value = '  -3.74%  '
$ws.Range('E13').Value = '  +0.98%  '
$ws.Range('D14').Value = '9.42'
$ws.Range('E14').Value = '  -2.18%  '
$ws.Range('D15').Value = '4.074.85'
$ws.Range('E15').Value = '  -1.53%  '
$ws.Range('D16').Value = '19.21'
$ws.Range('E16').Value = '  -3.21%  '
$ws.Range('D17').Value = '3.498.30'
$ws.Range('E17').Value = '  -2.08%  '
$ws.Range('D18').Value = '68.791.21'
$ws.Range('E18').Value = '  -1.23%  '
$ws.Range('D19').Value = '12.46'
$ws.Range('E19').Value = '  +0.36%  '
$ws.Range('D20').Value = '0.119'
$ws.Range('E20').Value = '  -1.13%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '534.46'
$ws.Range('E21').Value = '  +13.92%  '
$ws.Range('B22').Value = 'Polygon'
$ws.Range('C22').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D22').Value = '1.03'
$ws.Range('E22').Value = '  -0.62%  '
$ws.Range('D23').Value = '20.22'
$ws.Range('E23').Value = '  +4.52%  '
$ws.Range('D24').Value = '4.98'
$ws.Range('E24').Value = '  -1.75%  '
$ws.Range('E25').Value = '  +0.68%  '
$ws.Range('D26').Value = '93.82'
$ws.Range('E26').Value = '  +6.54%  '
$ws.Range('D27').Value = '10.99'
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E28').Value = '  -4.42%  '
$ws.Range('D29').Value = '9.12'
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('D30').Value = '31.60'
$ws.Range('E30').Value = '  -0.85%  '
$ws.Range('E31').Value = '  -4.31%  '
$ws.Range('D32').Value = '12.58'
$ws.Range('E32').Value = '  +4.64%  '
$ws.Range('D33').Value = '64.20'
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('D34').Value = '0.114'
$ws.Range('E34').Value = '  -4.48%  '
$ws.Range('D35').Value = '567.99'
$ws.Range('E35').Value = '  -2.57%  '
$ws.Range('D36').Value = '3.10'
$ws.Range('E36').Value = '  +8.53%  '
$ws.Range('E37').Value = '  -2.37%  '
$ws.Range('E38').Value = '  +0.09%  '
$ws.Range('E39').Value = '  +0.98%  '
$ws.Range('E40').Value = '  -4.97%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '0.133'
$ws.Range('E41').Value = '  -4.64%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '3.34'
$ws.Range('E42').Value = '  -4.62%  '
$ws.Range('D43').Value = '3.05'
$ws.Range('E43').Value = '  -4.75%  '
$ws.Range('B44').Value = 'ThetaToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D44').Value = '2.98'
$ws.Range('E44').Value = '  -3.88%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '3.48'
$ws.Range('E45').Value = '  +4.94%  '
$ws.Range('D46').Value = '0.0442'
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').Value = '3.161.03'
$ws.Range('E47').Value = '  -1.41%  '
$ws.Range('D48').Value = '9.16'
$ws.Range('E48').Value = '  -2.51%  '
$ws.Range('E49').Value = '  -1.99%  '
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('D51').Value = '137.13'
$ws.Range('E51').Value = '  -0.31%  '

# Restore the default (unstyled) appearance for column D now that the
# values are committed as text, so no stray style index lingers on cells.
$ws.Range("D2:D51").Style = "Normal"

"Updated cryptos list"